# Add data for 2021-11-10
# Updates the "through November 0X" report by one more day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab and update the header label in B1 to reflect the
# new "through" date.
$ws.Name = "Through 2021-11-02"
$ws.Range("B1").Value = "November 2021 (through November 02)"

# Row 2 - North Lawndale
$ws.Range("B2").Value = 1
$ws.Range("M2").Value = 2

# Row 3 - Garfield Park
$ws.Range("M3").Value = 3

# Row 4 - Austin
$ws.Range("AI4").Value = 2
$ws.Range("BP4").Value = 1

# Row 5 - Humboldt Park
$ws.Range("M5").Value = 1

# Row 7 - Englewood
$ws.Range("B7").Value = 1
$ws.Range("AI7").Value = 1

# Row 11 - New City
$ws.Range("B11").Value = 1
$ws.Range("AT11").Value = 1

# Row 16 - West Loop
$ws.Range("M16").Value = 1

# Row 17 - Calumet Heights
$ws.Range("B17").Value = 1

# Row 28 - Uptown
$ws.Range("X28").Value = 1
$ws.Range("AI28").Value = 2

# Row 29 - Near South Side
$ws.Range("AI29").Value = 1

# Row 41 - West Ridge
$ws.Range("B41").Value = 1

# Row 43 - Ashburn
$ws.Range("BE43").Value = 1

# Row 47 - Roseland
$ws.Range("X47").Value = 1
$ws.Range("AT47").Value = 1

# Row 65 - Brighton Park
$ws.Range("B65").Value = 1

# Row 84 - Morgan Park
$ws.Range("B84").Value = 2

# Row 97 - West Lawn
$ws.Range("BE97").Value = 1
